# Bugfixes and Thermocycler Improvements
# Update the last reagent row ("pAGM1299") with corrected concentration,
# fmol target, fragment/total base counts, and dilution factor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculator")
$ws.Activate()

$ws.Range("B8").Value = 308
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = 1646
$ws.Range("E8").Value = 2676
$ws.Range("H8").Value = 7

$ws.Range("H8").Select()
